# Update "Fresh bloom Flowers_2025-10-14.xlsx":
#  - Append 9 new order-line rows (42-50) plus a trailing summary row (51)
#    to the "Orders" sheet, extending the used range from A1:L41 to A1:L51.
#  - Append extra digits to the barcode-like string in "Summary"!G2.

$wb = $excel.ActiveWorkbook
$ordersWs = $wb.Worksheets.Item("Orders")
$summaryWs = $wb.Worksheets.Item("Summary")

# New order rows: row => (FlowerName, Number)
# The "Number" values are stored as text (matching the rest of the sheet),
# so a leading apostrophe forces text instead of Excel's auto numeric
# conversion for these digit-only strings.
$newRows = @(
    @{ Row = 42; Name = "238_苏菲宝贝_undefined_Rosa rugosa Thunb._10stems"; Number = "5" },
    @{ Row = 43; Name = "274_仙子之吻_undefined_Rosa rugosa Thunb._10stems"; Number = "7" },
    @{ Row = 44; Name = "268_猩红泡泡_spray red_Rosa rugosa Thunb._10stems"; Number = "5" },
    @{ Row = 45; Name = "13_酒红洋桔梗_Burgundy Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"; Number = "5" },
    @{ Row = 46; Name = "2_粉洋桔梗_Pink Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"; Number = "5" },
    @{ Row = 47; Name = "411_紫罗兰白_violet white_undefined_1bunch"; Number = "18" },
    @{ Row = 48; Name = "630_吸色康乃馨天蓝_tinted tiffany blue_undefined_20stems"; Number = "5" },
    @{ Row = 49; Name = "578_腊梅粉_wax pink_undefined_1bunch"; Number = "5" },
    @{ Row = 50; Name = "300_白星_White Gypso_ gypsophila_1kg"; Number = "15" }
)

foreach ($item in $newRows) {
    $ordersWs.Cells.Item($item.Row, 3).Value = $item.Name
    $ordersWs.Cells.Item($item.Row, 6).Value = "'" + $item.Number
}

# Trailing row 51: a lone package-count marker in column A.
$ordersWs.Cells.Item(51, 1).Value = "'6"

# Append extra digits to the barcode string on the Summary sheet.
$summaryWs.Cells.Item(2, 7).Value = "'058810310202055585103121565555517101015105551010101040104412575551855150"
